# Update the flag image paths from "Square/xx.png" to "World/xx.png"
# in the Flag column (C2:C22), and move the sheet selection to Q10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countryCodes = @(
    "is", "se", "fi", "gb", "tr", "es", "nl", "il", "ch", "ie",
    "de", "at", "dk", "gr", "no", "be", "lu", "it", "fr", "pt", "rs"
)

$startRow = 2
for ($i = 0; $i -lt $countryCodes.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = "World/" + $countryCodes[$i] + ".png"
}

# Update the selected cell shown in the saved sheet view.
$ws.Range("Q10").Select()
